$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns whose values participate in this record re-ordering
$cols = @("A", "B", "E", "F", "G", "H", "Q", "R", "AC")

# Swap rows 11 and 12 (records were swapped)
foreach ($col in $cols) {
    $addr11 = "${col}11"
    $addr12 = "${col}12"
    $v11 = $ws.Range($addr11).Value()
    $v12 = $ws.Range($addr12).Value()
    $ws.Range($addr11).Value = $v12
    $ws.Range($addr12).Value = $v11
}

# Rotate rows 17, 18, 19: new17 = old19, new18 = old17, new19 = old18
foreach ($col in $cols) {
    $addr17 = "${col}17"
    $addr18 = "${col}18"
    $addr19 = "${col}19"
    $v17 = $ws.Range($addr17).Value()
    $v18 = $ws.Range($addr18).Value()
    $v19 = $ws.Range($addr19).Value()
    $ws.Range($addr17).Value = $v19
    $ws.Range($addr18).Value = $v17
    $ws.Range($addr19).Value = $v18
}
